$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.822.54"
$ws.Range("E2").Value = "  -1.87%  "
$ws.Range("D3").Value = "2.907.88"
$ws.Range("E3").Value = "  -2.95%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'526.44"
$ws.Range("E5").Value = "  -3.15%  "
$ws.Range("D6").Value = "'144.06"
$ws.Range("E6").Value = "  -5.79%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").Value = "'0.546"
$ws.Range("E8").Value = "  -4.25%  "
$ws.Range("D9").Value = "2.913.54"
$ws.Range("E9").Value = "  -3.29%  "
$ws.Range("E10").Value = "  -5.24%  "
$ws.Range("D11").Value = "'6.09"
$ws.Range("E11").Value = "  -1.49%  "
$ws.Range("D12").Value = "'0.358"
$ws.Range("E12").Value = "  -3.15%  "
$ws.Range("D13").Value = "3.416.36"
$ws.Range("E13").Value = "  -3.04%  "
$ws.Range("E14").Value = "  +2.94%  "
$ws.Range("D15").Value = "60.768.46"
$ws.Range("E15").Value = "  -1.95%  "
$ws.Range("D16").Value = "'22.48"
$ws.Range("E16").Value = "  -6.43%  "
$ws.Range("D17").Value = "2.916.14"
$ws.Range("E17").Value = "  -2.90%  "
$ws.Range("D18").Value = "'0.0000141"
$ws.Range("E18").Value = "  -4.26%  "
$ws.Range("D19").Value = "'4.89"
$ws.Range("E19").Value = "  -5.55%  "
$ws.Range("D20").Value = "'11.57"
$ws.Range("E20").Value = "  -4.33%  "
$ws.Range("D21").Value = "'352.71"
$ws.Range("E21").Value = "  -7.56%  "
$ws.Range("D22").Value = "'6.50"
$ws.Range("E22").Value = "  -4.00%  "
$ws.Range("D23").Value = "'0.998"
$ws.Range("D24").Value = "'5.73"
$ws.Range("E24").Value = "  +1.23%  "
$ws.Range("D25").Value = "'64.81"
$ws.Range("E25").Value = "  -1.98%  "
$ws.Range("D26").Value = "'0.451"
$ws.Range("E26").Value = "  -4.25%  "
$ws.Range("E27").Value = "  -7.22%  "
$ws.Range("D28").Value = "'0.999"
$ws.Range("E28").Value = "  +0.21%  "
$ws.Range("D29").Value = "'7.84"
$ws.Range("E29").Value = "  -5.13%  "
$ws.Range("D30").Value = "0.0₃0866"
$ws.Range("E30").Value = "  -8.75%  "
$ws.Range("E31").Value = "  +0.03%  "
$ws.Range("E32").Value = "  -3.01%  "
$ws.Range("D33").Value = "'19.66"
$ws.Range("E33").Value = "  -4.38%  "
$ws.Range("D34").Value = "'152.33"
$ws.Range("E34").Value = "  -5.02%  "
$ws.Range("D35").Value = "'4.39"
$ws.Range("E35").Value = "  -5.44%  "
$ws.Range("E36").Value = "  -7.58%  "
$ws.Range("D37").Value = "'0.996"
$ws.Range("E37").Value = "  -7.70%  "
$ws.Range("D38").Value = "'1.20"
$ws.Range("E38").Value = "  -7.14%  "
$ws.Range("D39").Value = "'37.58"
$ws.Range("E39").Value = "  -0.40%  "
$ws.Range("D40").Value = "'1.47"
$ws.Range("E40").Value = "  -5.85%  "
$ws.Range("D41").Value = "'3.71"
$ws.Range("E41").Value = "  -5.45%  "
$ws.Range("D42").Value = "'0.652"
$ws.Range("E42").Value = "  -3.48%  "
$ws.Range("D43").Value = "2.278.01"
$ws.Range("E43").Value = "  -6.20%  "
$ws.Range("D44").Value = "'0.0581"
$ws.Range("E44").Value = "  -2.49%  "
$ws.Range("D45").Value = "'20.33"
$ws.Range("E45").Value = "  -8.82%  "
$ws.Range("E46").Value = "  +0.03%  "
$ws.Range("D47").Value = "'4.94"
$ws.Range("E47").Value = "  -5.20%  "
$ws.Range("D48").Value = "'0.0237"
$ws.Range("E48").Value = "  -4.11%  "
$ws.Range("E49").Value = "  -0.96%  "
$ws.Range("D50").Value = "'0.0918"
$ws.Range("E50").Value = "  -4.09%  "
$ws.Range("D51").Value = "'18.41"
$ws.Range("E51").Value = "  -7.84%  "
